$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price-column cells that look like plain numbers stay as text,
# matching the source data which stores prices as literal strings
# (e.g. "305.68", "10.10", "0.000008640") rather than numeric values.

$ws.Range('D2').Value = '26.398.72'
$ws.Range('E2').Value = '  -3.81%  '
$ws.Range('D3').Value = '1.766.82'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '305.68'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4303'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3636'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8454'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.27'
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D12').Value = '1.798.32'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.254'
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.434'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06803'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.24'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008640'
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.05'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Value = '26.401.77'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.039'
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.23'
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('D24').Value = '1.989.08'
$ws.Range('E24').Value = '  -4.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.65'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.855'
$ws.Range('E26').Value = '  -6.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.12'
$ws.Range('E27').Value = '  -3.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.079'
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.16'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08930'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7315'
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.114'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.333'
$ws.Range('E34').Value = '  -4.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.768'
$ws.Range('E35').Value = '  -6.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.073'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05126'
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01893'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1612'
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4915'
$ws.Range('E41').Value = '  -2.87%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.522'
$ws.Range('E42').Value = '  -9.26%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.232'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.047'
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.89'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.10'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06190'
$ws.Range('E48').Value = '  -3.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4477'
$ws.Range('E49').Value = '  -4.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.580'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('E51').Value = '  +2.07%  '
